$wb = $excel.ActiveWorkbook

# ---- Sheet: PUTWALL PICKING ----
$ws = $wb.Worksheets.Item("PUTWALL PICKING")
$ws.Range("A1:C12").ClearContents()
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "PutwallPickingQuantity"
$ws.Range("C1").Value = "UPH"
$ws.Range("A2").Value = "AGNE8120.CARUTH"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.3
$ws.Range("A3").Value = "AHME0710.JUBRAN"
$ws.Range("B3").Value = 79
$ws.Range("C3").Value = 23.7
$ws.Range("A4").Value = "BOHD0676.KUSHLIAK"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 3.6
$ws.Range("A5").Value = "DIAN4065.ENTRIALGO"
$ws.Range("B5").Value = 154
$ws.Range("C5").Value = 46.2
$ws.Range("A6").Value = "KADE3054.ZONGO"
$ws.Range("B6").Value = 38
$ws.Range("C6").Value = 11.4
$ws.Range("A7").Value = "LOANA.MBONGO"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 1.2
$ws.Range("A8").Value = "MAKEDA.OLLIVIERRE"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.3
$ws.Range("A9").Value = "PATI2298.ATSIANGBE"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 0.8999999999999999
$ws.Range("A10").Value = "SURESH.DHAWAN"
$ws.Range("B10").Value = 39
$ws.Range("C10").Value = 11.7
$ws.Range("A11").Value = "THIE6554.DIALLO"
$ws.Range("B11").Value = 56
$ws.Range("C11").Value = 16.8
$ws.Range("A12").Value = "XUAN754N.LU"
$ws.Range("B12").Value = 353
$ws.Range("C12").Value = 105.9

# ---- Sheet: REGULAR PICK ----
$ws = $wb.Worksheets.Item("REGULAR PICK")
$ws.Range("A1:C14").ClearContents()
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "RegularPickQuantity"
$ws.Range("C1").Value = "UPH"
$ws.Range("A2").Value = "ADOL798N.SEEMANNVAZQ"
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 12
$ws.Range("A3").Value = "AGNE8120.CARUTH"
$ws.Range("B3").Value = 33
$ws.Range("C3").Value = 9.9
$ws.Range("A4").Value = "BOHD0676.KUSHLIAK"
$ws.Range("B4").Value = 98
$ws.Range("C4").Value = 29.4
$ws.Range("A5").Value = "DIAN4065.ENTRIALGO"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 0.8999999999999999
$ws.Range("A6").Value = "IREN797N.CABRERA"
$ws.Range("B6").Value = 27
$ws.Range("C6").Value = 8.1
$ws.Range("A7").Value = "JEEW9554.SITUMUDALIG"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 3.3
$ws.Range("A8").Value = "LOANA.MBONGO"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 1.2
$ws.Range("A9").Value = "MAKEDA.OLLIVIERRE"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 1.5
$ws.Range("A10").Value = "PATI2298.ATSIANGBE"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 0.8999999999999999
$ws.Range("A11").Value = "PATR5027.AMEH"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 0.8999999999999999
$ws.Range("A12").Value = "WESL5337.CADETTE"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 0.8999999999999999
$ws.Range("A13").Value = "XUAN754N.LU"
$ws.Range("B13").Value = 15
$ws.Range("C13").Value = 4.5
$ws.Range("A14").Value = "ZAKI0190.PHILLIPHORS"
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 3.6

# ---- Sheet: SINGLE PICK ----
$ws = $wb.Worksheets.Item("SINGLE PICK")
$ws.Range("A1:C10").ClearContents()
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "SinglePickQuantity"
$ws.Range("C1").Value = "UPH"
$ws.Range("A2").Value = "BUDD0680.TENNAKOON"
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = 8.7
$ws.Range("A3").Value = "JEEW9554.SITUMUDALIG"
$ws.Range("B3").Value = 27
$ws.Range("C3").Value = 8.1
$ws.Range("A4").Value = "LOWRHY-OTIENO.JAOKO"
$ws.Range("B4").Value = 90
$ws.Range("C4").Value = 27
$ws.Range("A5").Value = "OMAR6689.KHAN"
$ws.Range("B5").Value = 69
$ws.Range("C5").Value = 20.7
$ws.Range("A6").Value = "PATR5027.AMEH"
$ws.Range("B6").Value = 60
$ws.Range("C6").Value = 18
$ws.Range("A7").Value = "PRINCE.FORSON"
$ws.Range("B7").Value = 18
$ws.Range("C7").Value = 5.399999999999999
$ws.Range("A8").Value = "WESL5337.CADETTE"
$ws.Range("B8").Value = 90
$ws.Range("C8").Value = 27
$ws.Range("A9").Value = "WILDINE.JEUNE"
$ws.Range("B9").Value = 80
$ws.Range("C9").Value = 24

# ---- Sheet: REPLENISHMENT PICK ----
$ws = $wb.Worksheets.Item("REPLENISHMENT PICK")
$ws.Range("A1:C21").ClearContents()
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "ReplenishmentPickQuantity"
$ws.Range("C1").Value = "UPH"
$ws.Range("A2").Value = "AHME0710.JUBRAN"
$ws.Range("B2").Value = 31
$ws.Range("C2").Value = 9.299999999999999
$ws.Range("A3").Value = "BOHD0676.KUSHLIAK"
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 3.3
$ws.Range("A4").Value = "BUDD0680.TENNAKOON"
$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 12
$ws.Range("A5").Value = "DEVI789.SINGH"
$ws.Range("B5").Value = 64
$ws.Range("C5").Value = 19.2
$ws.Range("A6").Value = "DIAN4065.ENTRIALGO"
$ws.Range("B6").Value = 119
$ws.Range("C6").Value = 35.7
$ws.Range("A7").Value = "IREN797N.CABRERA"
$ws.Range("B7").Value = 68
$ws.Range("C7").Value = 20.4
$ws.Range("A8").Value = "JEEW9554.SITUMUDALIG"
$ws.Range("B8").Value = 92
$ws.Range("C8").Value = 27.6
$ws.Range("A9").Value = "JIGN776N.PATEL"
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 7.5
$ws.Range("A10").Value = "KADE3054.ZONGO"
$ws.Range("B10").Value = 71
$ws.Range("C10").Value = 21.3
$ws.Range("A11").Value = "LOWRHY-OTIENO.JAOKO"
$ws.Range("B11").Value = 39
$ws.Range("C11").Value = 11.7
$ws.Range("A12").Value = "MICA0432.RIZKALLAMAR"
$ws.Range("B12").Value = 171
$ws.Range("C12").Value = 51.3
$ws.Range("A13").Value = "NESR2403.ATTALAH"
$ws.Range("B13").Value = 23
$ws.Range("C13").Value = 6.899999999999999
$ws.Range("A14").Value = "OMAR6689.KHAN"
$ws.Range("B14").Value = 42
$ws.Range("C14").Value = 12.6
$ws.Range("A15").Value = "PATR5027.AMEH"
$ws.Range("B15").Value = 45
$ws.Range("C15").Value = 13.5
$ws.Range("A16").Value = "PRINCE.FORSON"
$ws.Range("B16").Value = 35
$ws.Range("C16").Value = 10.5
$ws.Range("A17").Value = "THIE6554.DIALLO"
$ws.Range("B17").Value = 94
$ws.Range("C17").Value = 28.2
$ws.Range("A18").Value = "WESL5337.CADETTE"
$ws.Range("B18").Value = 38
$ws.Range("C18").Value = 11.4
$ws.Range("A19").Value = "WILDINE.JEUNE"
$ws.Range("B19").Value = 40
$ws.Range("C19").Value = 12
$ws.Range("A20").Value = "YATI0689.YATIN"
$ws.Range("B20").Value = 77
$ws.Range("C20").Value = 23.1
$ws.Range("A21").Value = "ZAKI0190.PHILLIPHORS"
$ws.Range("B21").Value = 149
$ws.Range("C21").Value = 44.7

# ---- Sheet: IDLE TIME ----
$ws = $wb.Worksheets.Item("IDLE TIME")
$ws.Range("A1:B28").ClearContents()
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "TotalIdleTime"
$ws.Range("A2").Value = "ADOL798N.SEEMANNVAZQ"
$ws.Range("B2").Value = 159
$ws.Range("A3").Value = "AGNE8120.CARUTH"
$ws.Range("B3").Value = 181
$ws.Range("A4").Value = "AHME0710.JUBRAN"
$ws.Range("B4").Value = 151
$ws.Range("A5").Value = "BOHD0676.KUSHLIAK"
$ws.Range("B5").Value = 63
$ws.Range("A6").Value = "BUDD0680.TENNAKOON"
$ws.Range("B6").Value = 114
$ws.Range("A7").Value = "DEVI789.SINGH"
$ws.Range("B7").Value = 116
$ws.Range("A8").Value = "DIAN4065.ENTRIALGO"
$ws.Range("B8").Value = 54
$ws.Range("A9").Value = "IREN797N.CABRERA"
$ws.Range("B9").Value = 80
$ws.Range("A10").Value = "JEEW9554.SITUMUDALIG"
$ws.Range("B10").Value = 71
$ws.Range("A11").Value = "JIGN776N.PATEL"
$ws.Range("B11").Value = 122
$ws.Range("A12").Value = "KADE3054.ZONGO"
$ws.Range("B12").Value = 101
$ws.Range("A13").Value = "LOANA.MBONGO"
$ws.Range("B13").Value = 196
$ws.Range("A14").Value = "LOWRHY-OTIENO.JAOKO"
$ws.Range("B14").Value = 53
$ws.Range("A15").Value = "MAKEDA.OLLIVIERRE"
$ws.Range("B15").Value = 185
$ws.Range("A16").Value = "MICA0432.RIZKALLAMAR"
$ws.Range("B16").Value = 64
$ws.Range("A17").Value = "NESR2403.ATTALAH"
$ws.Range("B17").Value = 191
$ws.Range("A18").Value = "OMAR6689.KHAN"
$ws.Range("B18").Value = 51
$ws.Range("A19").Value = "PATI2298.ATSIANGBE"
$ws.Range("B19").Value = 134
$ws.Range("A20").Value = "PATR5027.AMEH"
$ws.Range("B20").Value = 65
$ws.Range("A21").Value = "PRINCE.FORSON"
$ws.Range("B21").Value = 188
$ws.Range("A22").Value = "SURESH.DHAWAN"
$ws.Range("B22").Value = 132
$ws.Range("A23").Value = "THIE6554.DIALLO"
$ws.Range("B23").Value = 99
$ws.Range("A24").Value = "WESL5337.CADETTE"
$ws.Range("B24").Value = 102
$ws.Range("A25").Value = "WILDINE.JEUNE"
$ws.Range("B25").Value = 51
$ws.Range("A26").Value = "XUAN754N.LU"
$ws.Range("B26").Value = 44
$ws.Range("A27").Value = "YATI0689.YATIN"
$ws.Range("B27").Value = 135
$ws.Range("A28").Value = "ZAKI0190.PHILLIPHORS"
$ws.Range("B28").Value = 49

# ---- Sheet: Total Units picked by hour ----
$ws = $wb.Worksheets.Item("Total Units picked by hour")
$ws.Range("A1:E6").ClearContents()
$ws.Range("A1").Value = "Hour"
$ws.Range("B1").Value = "Regular Pick"
$ws.Range("C1").Value = "Single Pick"
$ws.Range("D1").Value = "Replenishment Pick"
$ws.Range("E1").Value = "Putwall Pick"
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = -18
$ws.Range("C2").Value = -5
$ws.Range("D2").Value = -313
$ws.Range("E2").Value = -119
$ws.Range("A3").Value = 21
$ws.Range("B3").Value = -78
$ws.Range("C3").Value = -111
$ws.Range("D3").Value = -413
$ws.Range("E3").Value = -18
$ws.Range("A4").Value = 22
$ws.Range("B4").Value = -99
$ws.Range("C4").Value = -232
$ws.Range("D4").Value = -373
$ws.Range("E4").Value = -356
$ws.Range("A5").Value = 23
$ws.Range("B5").Value = -62
$ws.Range("C5").Value = -115
$ws.Range("D5").Value = -175
$ws.Range("E5").Value = -247
$ws.Range("A6").Value = "Total"
$ws.Range("B6").Value = -257
$ws.Range("C6").Value = -463
$ws.Range("D6").Value = -1274
$ws.Range("E6").Value = -740
